$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.117.33'
$ws.Range("E2").Value = '  +2.89%  '
$ws.Range("D3").Value = '2.651.99'
$ws.Range("E3").Value = '  +2.72%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.76%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.118'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.81%  '
$ws.Range("E10").Value = '  +4.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.79'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.03'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000188'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +21.48%  '
$ws.Range("D15").Value = '3.125.22'
$ws.Range("E15").Value = '  +2.57%  '
$ws.Range("D16").Value = '65.009.93'
$ws.Range("E16").Value = '  +3.04%  '
$ws.Range("D17").Value = '2.664.48'
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.81'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.41%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.33%  '
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.52'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.18%  '
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("E27").Value = '  +1.43%  '
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("D29").Value = '0.0₃0950'
$ws.Range("E29").Value = '  +10.53%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '522.81'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.35%  '
$ws.Range("E32").Value = '  +3.82%  '
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.63'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.45%  '
$ws.Range("E35").Value = '  +4.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.428'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '165.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '165.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("E44").Value = '  +2.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0619'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.54%  '
$ws.Range("E47").Value = '  +4.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.649'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.52%  '
$ws.Range("E49").Value = '  +1.58%  '
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.78%  '
